$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set new header values P1/Q1 (will be styled below to match O1)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("B2").Value = 3.46677384630982
$ws.Range("C2").Value = 1.019169272890139
$ws.Range("D2").Value = 0.04824190763466873
$ws.Range("E2").Value = 1.35169962688731
$ws.Range("F2").Value = 0.5084535737087279
$ws.Range("G2").Value = 0.0007896385517037115
$ws.Range("H2").Value = 0.01092094002687105
$ws.Range("I2").Value = 0.003536527785124033
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 1.435151477765999
$ws.Range("B3").Value = 3.02011174815118
$ws.Range("C3").Value = 0.8990605546520669
$ws.Range("D3").Value = 0.0433835340809452
$ws.Range("E3").Value = 1.176740016824183
$ws.Range("F3").Value = 0.468187038062311
$ws.Range("G3").Value = 0.000793366441376707
$ws.Range("H3").Value = 0.00766975902495437
$ws.Range("I3").Value = 0.001965288275290966
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 1.348594257482659
$ws.Range("B4").Value = 2.745669458445775
$ws.Range("C4").Value = 0.8258832939075091
$ws.Range("D4").Value = 0.04041041276209967
$ws.Range("E4").Value = 1.069666178277558
$ws.Range("F4").Value = 0.4441198563540496
$ws.Range("G4").Value = 0.0007957272508090085
$ws.Range("H4").Value = 0.005904022270549425
$ws.Range("I4").Value = 0.001270814676121557
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 1.297464469391542
$ws.Range("B5").Value = 2.633730900298815
$ws.Range("C5").Value = 0.7977722387358028
$ws.Range("D5").Value = 0.03928077046956702
$ws.Range("E5").Value = 1.02608314129057
$ws.Range("F5").Value = 0.4336700664951891
$ws.Range("G5").Value = 0.0007967124296093884
$ws.Range("H5").Value = 0.00523637476029204
$ws.Range("I5").Value = 0.001110499418202515
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 1.274436995028537
$ws.Range("B6").Value = 2.615103766652396
$ws.Range("C6").Value = 0.7950489896263662
$ws.Range("D6").Value = 0.03919122151063448
$ws.Range("E6").Value = 1.018830727099584
$ws.Range("F6").Value = 0.4309748955527013
$ws.Range("G6").Value = 0.0007968828780179578
$ws.Range("H6").Value = 0.005126972417599962
$ws.Range("I6").Value = 0.001168165231651308
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 1.267397047251904
$ws.Range("B7").Value = 2.74406360535454
$ws.Range("C7").Value = 0.8308092758339853
$ws.Range("D7").Value = 0.04066341026042863
$ws.Range("E7").Value = 1.069026406947373
$ws.Range("F7").Value = 0.4413300423773592
$ws.Range("G7").Value = 0.0007957557770997617
$ws.Range("H7").Value = 0.005890106154369357
$ws.Range("I7").Value = 0.001472853056091239
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 1.28828808501018
$ws.Range("B8").Value = 3.312653545137607
$ws.Range("C8").Value = 0.9847298999007421
$ws.Range("D8").Value = 0.04692401796475565
$ws.Range("E8").Value = 1.291216389130142
$ws.Range("F8").Value = 0.4909136236972671
$ws.Range("G8").Value = 0.0007909279253908628
$ws.Range("H8").Value = 0.009741767389663553
$ws.Range("I8").Value = 0.003174773938467546
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 1.393125671300481
$ws.Range("B9").Value = 4.427610317638539
$ws.Range("C9").Value = 1.282110248089054
$ws.Range("D9").Value = 0.0588090931990024
$ws.Range("E9").Value = 1.730915313986756
$ws.Range("F9").Value = 0.5988199823157387
$ws.Range("G9").Value = 0.0007819898167375172
$ws.Range("H9").Value = 0.01935802986678037
$ws.Range("I9").Value = 0.008607690668569923
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 1.632718893241275
$ws.Range("B10").Value = 5.24749449069725
$ws.Range("C10").Value = 1.50644471872306
$ws.Range("D10").Value = 0.06925275944912102
$ws.Range("E10").Value = 1.955498156088424
$ws.Range("F10").Value = 0.6683724546327312
$ws.Range("G10").Value = 0.0007759177289413412
$ws.Range("H10").Value = 0.02725836382229918
$ws.Range("I10").Value = 0.0143664161419883
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 1.775547518866148
$ws.Range("B11").Value = 5.612142301903873
$ws.Range("C11").Value = 1.610434802784823
$ws.Range("D11").Value = 0.08861588644301577
$ws.Range("E11").Value = 1.267337406619561
$ws.Range("F11").Value = 0.5871507068896307
$ws.Range("G11").Value = 0.0007746196692268896
$ws.Range("H11").Value = 0.04270868006982909
$ws.Range("I11").Value = 0.01643561128637661
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 1.474547141099578
$ws.Range("B12").Value = 5.74725555498668
$ws.Range("C12").Value = 1.642389293035194
$ws.Range("D12").Value = 0.1037805964393641
$ws.Range("E12").Value = 0.7713189373199754
$ws.Range("F12").Value = 0.5121260620343264
$ws.Range("G12").Value = 0.0007745905169331831
$ws.Range("H12").Value = 0.07854308646568597
$ws.Range("I12").Value = 0.01664005842228544
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 1.224683832186997
$ws.Range("B13").Value = 5.71206086016997
$ws.Range("C13").Value = 1.628934670659646
$ws.Range("D13").Value = 0.1168170525753141
$ws.Range("E13").Value = 0.3948377359290731
$ws.Range("F13").Value = 0.4334589186537912
$ws.Range("G13").Value = 0.0007755528712741501
$ws.Range("H13").Value = 0.1314041879758321
$ws.Range("I13").Value = 0.01567717074598907
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 0.9858811040282944
$ws.Range("B14").Value = 5.610402796626545
$ws.Range("C14").Value = 1.600782500645209
$ws.Range("D14").Value = 0.1251548963565625
$ws.Range("E14").Value = 0.2040035724406977
$ws.Range("F14").Value = 0.3772055239524477
$ws.Range("G14").Value = 0.0007766576504858142
$ws.Range("H14").Value = 0.1788295655706946
$ws.Range("I14").Value = 0.01459786682350295
$ws.Range("O14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 0.8260307867355579
$ws.Range("B15").Value = 5.549222564864976
$ws.Range("C15").Value = 1.586479300100507
$ws.Range("D15").Value = 0.1266680095887693
$ws.Range("E15").Value = 0.1664673984577796
$ws.Range("F15").Value = 0.361549071281118
$ws.Range("G15").Value = 0.0007771664842525621
$ws.Range("H15").Value = 0.1906857109769646
$ws.Range("I15").Value = 0.01418120435697556
$ws.Range("O15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 0.7848470606577962
$ws.Range("B16").Value = 5.202401853822494
$ws.Range("C16").Value = 1.495615476278033
$ws.Range("D16").Value = 0.1191748546277722
$ws.Range("E16").Value = 0.1619581450697751
$ws.Range("F16").Value = 0.3473887037929586
$ws.Range("G16").Value = 0.0007795445118418047
$ws.Range("H16").Value = 0.175491235802042
$ws.Range("I16").Value = 0.01192358401845794
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 0.7706965671103063
$ws.Range("B17").Value = 4.991035451832488
$ws.Range("C17").Value = 1.441267988981963
$ws.Range("D17").Value = 0.1089952125854552
$ws.Range("E17").Value = 0.2459753982704243
$ws.Range("F17").Value = 0.365796307107324
$ws.Range("G17").Value = 0.000780802229109395
$ws.Range("H17").Value = 0.1367689001963726
$ws.Range("I17").Value = 0.01081278822335374
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 0.8422055574621652
$ws.Range("B18").Value = 4.871728470037112
$ws.Range("C18").Value = 1.407347280808153
$ws.Range("D18").Value = 0.09563455489917061
$ws.Range("E18").Value = 0.4790218023010837
$ws.Range("F18").Value = 0.4181691898324047
$ws.Range("G18").Value = 0.0007811315644281613
$ws.Range("H18").Value = 0.08461431799794639
$ws.Range("I18").Value = 0.01026415122278657
$ws.Range("O18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 1.011170216096644
$ws.Range("B19").Value = 4.83500880465499
$ws.Range("C19").Value = 1.40339474351282
$ws.Range("D19").Value = 0.08223365033970964
$ws.Range("E19").Value = 0.9131096264680565
$ws.Range("F19").Value = 0.4940337535870114
$ws.Range("G19").Value = 0.0007805737522911431
$ws.Range("H19").Value = 0.0420758044946794
$ws.Range("I19").Value = 0.01072812603060136
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 1.251779364725252
$ws.Range("B20").Value = 5.031627643970751
$ws.Range("C20").Value = 1.464437765624723
$ws.Range("D20").Value = 0.06744620694644254
$ws.Range("E20").Value = 1.891890550245648
$ws.Range("F20").Value = 0.6409981226065185
$ws.Range("G20").Value = 0.0007775358304290444
$ws.Range("H20").Value = 0.02500230937025849
$ws.Range("I20").Value = 0.01333659226125761
$ws.Range("O20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 1.707873429777266
$ws.Range("B21").Value = 5.663603663768583
$ws.Range("C21").Value = 1.636139845533137
$ws.Range("D21").Value = 0.07300782904698622
$ws.Range("E21").Value = 2.224275313971731
$ws.Range("F21").Value = 0.7181090159483858
$ws.Range("G21").Value = 0.0007726865873697398
$ws.Range("H21").Value = 0.03281830044685519
$ws.Range("I21").Value = 0.01854060229383236
$ws.Range("O21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 1.895197017569245
$ws.Range("B22").Value = 6.077296100170315
$ws.Range("C22").Value = 1.741715857843474
$ws.Range("D22").Value = 0.07705961331203781
$ws.Range("E22").Value = 2.390794747381719
$ws.Range("F22").Value = 0.7665380570823004
$ws.Range("G22").Value = 0.000769642109448282
$ws.Range("H22").Value = 0.03795447932891705
$ws.Range("I22").Value = 0.02207615676014552
$ws.Range("O22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 2.010982655735347
$ws.Range("B23").Value = 5.856533976777541
$ws.Range("C23").Value = 1.678693843975964
$ws.Range("D23").Value = 0.0745603636229859
$ws.Range("E23").Value = 2.301860511110718
$ws.Range("F23").Value = 0.7437029678328884
$ws.Range("G23").Value = 0.0007712498272074986
$ws.Range("H23").Value = 0.03518443906473134
$ws.Range("I23").Value = 0.01990882957891671
$ws.Range("O23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 1.959238721746175
$ws.Range("B24").Value = 5.022345505370311
$ws.Range("C24").Value = 1.452387582760707
$ws.Range("D24").Value = 0.0656626206233426
$ws.Range("E24").Value = 1.967538615152051
$ws.Range("F24").Value = 0.6548386594035378
$ws.Range("G24").Value = 0.0007774558731151467
$ws.Range("H24").Value = 0.02548343545092102
$ws.Range("I24").Value = 0.01297716825878936
$ws.Range("O24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 1.754635322995938
$ws.Range("B25").Value = 4.125615602166818
$ws.Range("C25").Value = 1.210872665251486
$ws.Range("D25").Value = 0.05608168448596018
$ws.Range("E25").Value = 1.611287029704712
$ws.Range("F25").Value = 0.5640920550256467
$ws.Range("G25").Value = 0.0007843738772956623
$ws.Range("H25").Value = 0.01649207729639268
$ws.Range("I25").Value = 0.007186637996994705
$ws.Range("O25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 1.549647283128024

# Apply header style (bold/centered/bordered, matching O1) to new header cells P1:Q1
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
